$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 899.8
$ws.Range("I4").Value = 1074.75
$ws.Range("K4").Value = 1074.75
$ws.Range("M4").Value = -960.75
$ws.Range("H17").Value = 1727.0
$ws.Range("J17").Value = 1727.0
$ws.Range("L17").Value = 5181.0
$ws.Range("N17").Value = -5517.0
$ws.Range("H51").Value = 77733.0
$ws.Range("I51").Value = 3200.0
$ws.Range("K51").Value = 3200.0
$ws.Range("M51").Value = -2716.0
$ws.Range("H62").Value = 11675.941
$ws.Range("J62").Value = 12963.286
$ws.Range("L62").Value = 12963.286
$ws.Range("N62").Value = -14211.286
$ws.Range("H65").Value = 11675.941
$ws.Range("J65").Value = 12963.286
$ws.Range("L65").Value = 64816.43
$ws.Range("N65").Value = -71056.43

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 455.75
$ws.Range("I5").Value = 346.625
$ws.Range("J5").Value = 674.0
$ws.Range("K5").Value = 346.625
$ws.Range("L5").Value = 674.0
$ws.Range("M5").Value = -234.625
$ws.Range("N5").Value = -898.0
$ws.Range("H31").Value = 4528.6665
$ws.Range("J31").Value = 0.0
$ws.Range("L31").Value = 0.0
$ws.Range("N31").ClearContents()
$ws.Range("H32").Value = 24100.04
$ws.Range("J32").Value = 67202.7
$ws.Range("L32").Value = 67202.7
$ws.Range("N32").Value = -67776.7
$ws.Range("H74").Value = 7498.189
$ws.Range("I74").Value = 2917.2222
$ws.Range("J74").Value = 19866.8
$ws.Range("K74").Value = 2917.2222
$ws.Range("L74").Value = 19866.8
$ws.Range("M74").Value = -2043.2222
$ws.Range("N74").Value = -21614.8
$ws.Range("H77").Value = 7498.189
$ws.Range("I77").Value = 2917.2222
$ws.Range("J77").Value = 19866.8
$ws.Range("K77").Value = 14586.111
$ws.Range("L77").Value = 99334.0
$ws.Range("M77").Value = -10218.111
$ws.Range("N77").Value = -108070.0
$ws.Range("H94").Value = 33219.832
$ws.Range("J94").Value = 33219.832
$ws.Range("L94").Value = 33219.832
$ws.Range("N94").Value = -35021.832
$ws.Range("H97").Value = 1299.0
$ws.Range("I97").Value = 1299.0
$ws.Range("J97").Value = 0.0
$ws.Range("K97").Value = 1299.0
$ws.Range("L97").Value = 0.0
$ws.Range("M97").Value = -803.0
$ws.Range("N97").ClearContents()
$ws.Range("H106").Value = 52728.75
$ws.Range("J106").Value = 52728.75
$ws.Range("L106").Value = 52728.75
$ws.Range("N106").Value = -55252.75
$ws.Range("H119").Value = 60886.11
$ws.Range("J119").Value = 60886.11
$ws.Range("L119").Value = 60886.11
$ws.Range("N119").Value = -70562.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 455.75
$ws.Range("I4").Value = 346.625
$ws.Range("J4").Value = 674.0
$ws.Range("K4").Value = 346.625
$ws.Range("L4").Value = 674.0
$ws.Range("M4").Value = -231.625
$ws.Range("N4").Value = -904.0
$ws.Range("H15").Value = 25000.0
$ws.Range("J15").Value = 25000.0
$ws.Range("L15").Value = 25000.0
$ws.Range("N15").Value = -25454.0
$ws.Range("H19").Value = 24999.5
$ws.Range("J19").Value = 24999.5
$ws.Range("L19").Value = 24999.5
$ws.Range("N19").Value = -25345.5
$ws.Range("H99").Value = 6791.1665
$ws.Range("I99").Value = 9591.077
$ws.Range("J99").Value = 3482.182
$ws.Range("K99").Value = 9591.077
$ws.Range("L99").Value = 3482.182
$ws.Range("M99").Value = -8093.076999999999
$ws.Range("N99").Value = -6478.182
$ws.Range("H107").Value = 1943.7241
$ws.Range("I107").Value = 1532.909
$ws.Range("K107").Value = 1532.909
$ws.Range("M107").Value = 387.0909999999999
$ws.Range("H130").Value = 99979.0
$ws.Range("J130").Value = 99979.0
$ws.Range("L130").Value = 99979.0
$ws.Range("N130").Value = -110019.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 2000.0
$ws.Range("I23").Value = 2000.0
$ws.Range("K23").Value = 2000.0
$ws.Range("M23").Value = -1760.0
$ws.Range("H27").Value = 2000.0
$ws.Range("I27").Value = 2000.0
$ws.Range("K27").Value = 2000.0
$ws.Range("M27").Value = -1808.0
$ws.Range("H68").Value = 73790.0
$ws.Range("J68").Value = 73790.0
$ws.Range("L68").Value = 73790.0
$ws.Range("N68").Value = -75288.0
$ws.Range("H71").Value = 73790.0
$ws.Range("J71").Value = 73790.0
$ws.Range("L71").Value = 221370.0
$ws.Range("N71").Value = -228858.0
$ws.Range("H86").Value = 53217.332
$ws.Range("I86").Value = 5419.1
$ws.Range("K86").Value = 5419.1
$ws.Range("M86").Value = -4296.1
$ws.Range("H89").Value = 53217.332
$ws.Range("I89").Value = 5419.1
$ws.Range("K89").Value = 27095.5
$ws.Range("M89").Value = -21479.5
$ws.Range("H92").Value = 110500.0
$ws.Range("J92").Value = 110500.0
$ws.Range("L92").Value = 110500.0
$ws.Range("N92").Value = -115492.0
$ws.Range("H132").Value = 7014.3
$ws.Range("I132").Value = 3091.2
$ws.Range("K132").Value = 9273.599999999999
$ws.Range("M132").Value = -6743.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14375.429
$ws.Range("I131").Value = 9126.667
$ws.Range("K131").Value = 27380.001
$ws.Range("M131").Value = -22340.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J11").Value = 1258746.8
$ws.Range("L11").Value = 1258746.8
$ws.Range("N11").Value = -1259024.8
$ws.Range("H18").Value = 59995.0
$ws.Range("I18").Value = 59990.0
$ws.Range("K18").Value = 59990.0
$ws.Range("M18").Value = -59697.0
$ws.Range("H126").Value = 4550.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 541.61536
$ws.Range("I16").Value = 541.61536
$ws.Range("K16").Value = 541.61536
$ws.Range("M16").Value = -371.61536
$ws.Range("H46").Value = 3679.9524
$ws.Range("I46").Value = 2767.182
$ws.Range("J46").Value = 4684.0
$ws.Range("K46").Value = 2767.182
$ws.Range("L46").Value = 4684.0
$ws.Range("M46").Value = -2579.182
$ws.Range("N46").Value = -5060.0
$ws.Range("H68").Value = 3575.2
$ws.Range("I68").Value = 3575.2
$ws.Range("K68").Value = 3575.2
$ws.Range("M68").Value = -2826.2
$ws.Range("H71").Value = 3575.2
$ws.Range("I71").Value = 3575.2
$ws.Range("K71").Value = 17876.0
$ws.Range("M71").Value = -14132.0
$ws.Range("H82").Value = 2548.2144
$ws.Range("I82").Value = 1963.3334
$ws.Range("J82").Value = 2986.875
$ws.Range("K82").Value = 1963.3334
$ws.Range("L82").Value = 2986.875
$ws.Range("M82").Value = -1602.3334
$ws.Range("N82").Value = -3708.875
$ws.Range("H85").Value = 2548.2144
$ws.Range("I85").Value = 1963.3334
$ws.Range("J85").Value = 2986.875
$ws.Range("K85").Value = 1963.3334
$ws.Range("L85").Value = 2986.875
$ws.Range("M85").Value = -715.3334
$ws.Range("N85").Value = -5482.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0.0
$ws.Range("I15").Value = 0.0
$ws.Range("K15").Value = 0.0
$ws.Range("M15").ClearContents()
$ws.Range("H39").Value = 0.0
$ws.Range("I39").Value = 0.0
$ws.Range("K39").Value = 0.0
$ws.Range("M39").ClearContents()
$ws.Range("H42").Value = 0.0
$ws.Range("I42").Value = 0.0
$ws.Range("K42").Value = 0.0
$ws.Range("M42").ClearContents()
$ws.Range("H49").Value = 34694.4
$ws.Range("I49").Value = 28242.5
$ws.Range("K49").Value = 28242.5
$ws.Range("M49").Value = -28012.5
$ws.Range("H126").Value = 3502.0
$ws.Range("I126").Value = 3502.0
$ws.Range("K126").Value = 10506.0
$ws.Range("M126").Value = -8036.0
$ws.Range("H132").Value = 4907208.0
$ws.Range("I132").Value = 5807.357
$ws.Range("J132").Value = 27780412.0
$ws.Range("K132").Value = 17422.071
$ws.Range("L132").Value = 83341236.0
$ws.Range("M132").Value = -14892.071
$ws.Range("N132").Value = -83346296.0
$ws.Range("H136").Value = 4382.15
$ws.Range("I136").Value = 2582.4856
$ws.Range("J136").Value = 16979.8
$ws.Range("K136").Value = 7747.4568
$ws.Range("L136").Value = 50939.39999999999
$ws.Range("M136").Value = -5197.4568
$ws.Range("N136").Value = -56039.39999999999

